$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the NPD (non-taxable income) formula in column P for rows 5-14:
# threshold 350 -> 400, rate 0.17 -> 0.19 (effective 2020-07-01)
for ($r = 5; $r -le 14; $r++) {
    $old = $ws.Range("P$r").Formula
    $new = $old -replace '(?<![0-9.])350(?![0-9])', '400' -replace '(?<![0-9])0\.17(?![0-9])', '0.19'
    $ws.Range("P$r").Formula = $new
}

# Update the active selection to match the author's last selection
$ws.Range("G24").Select()
